# Insert a new data row at row 17 (pushing existing rows 17-155 down to 18-156)
# and populate it with the newly reported price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("17:17").Insert()

$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44670
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 100114007
$ws.Range("G17").Value = "Jengibre"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 15
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 20000
$ws.Range("N17").Value = "`$/caja 13 kilos"
$ws.Range("O17").Value = "Perú"
$ws.Range("P17").Value = 1538
$ws.Range("Q17").Value = 13
$ws.Range("R17").Value = "Hortaliza"
